$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row right after the current last data row.
$lastRow = $ws.Cells(($ws.Rows.Count), 1).End(-4162).Row
$newRow = $lastRow + 1

# Force the new row to be stored as plain text (same "numberStoredAsText"
# shape as the rest of the sheet) so date-/number-looking values such as
# "2025-11-17" or "251117" are not silently coerced into date/number cells.
$ws.Range("A" + $newRow + ":E" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-11-17"
$ws.Range("B" + $newRow).Value = "Pick 3"
$ws.Range("C" + $newRow).Value = "251117"
$ws.Range("D" + $newRow).Value = "3-5-6"
$ws.Range("E" + $newRow).Value = "2025-11-17T21:40:37.883+04:00"
